$wb = $excel.ActiveWorkbook

# --- Update the previously-last sheet (CreateNewCampaign): it was the
# selected tab before; after adding the new sheet it no longer is, and its
# selection moves to A2:B2. Selecting it now (before the new sheet is
# created) reproduces that, since adding+naming the new sheet afterwards
# makes the new sheet the active tab instead.
$campaignSheet = $wb.Worksheets.Item("CreateNewCampaign")
[void]$campaignSheet.Range("A2:B2").Select()

# --- Add the new "CreateNewEvent" sheet after the last existing sheet ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add([System.Type]::Missing, $lastSheet)
$newSheet.Name = "CreateNewEvent"

# Hyperlink on A2 (mirrors the mailto: link used on every other sheet).
# Must be added before the row is formatted/filled so the cell picks up
# the shared "Hyperlink" style (xf index 3) instead of creating a stray one.
$newSheet.Hyperlinks.Add($newSheet.Range("A2"), "mailto:adilkhaleque429@gmail.com")

# Force text storage (matches numFmtId 49 "@" used throughout the workbook)
# so numeric-looking values ("4", "1", "0") are written as shared strings
# instead of numbers.
$newSheet.Range("A2:P2").NumberFormat = "@"

$newSheet.Range("A2").Value = "adilkhaleque429@gmail.com"
$newSheet.Range("B2").Value = "Testunbound6F"
$newSheet.Range("C2").Value = "Event1"
$newSheet.Range("D2").Value = "02/01/2022 08:00"
$newSheet.Range("E2").Value = "02/01/2022 08:30"
$newSheet.Range("F2").Value = "4"
$newSheet.Range("G2").Value = "This event includes Deal1, Task1, and Case1."
$newSheet.Range("H2").Value = "New York, NY"
$newSheet.Range("I2").Value = "Deal1"
$newSheet.Range("J2").Value = "Task1"
$newSheet.Range("K2").Value = "Case1"
$newSheet.Range("L2").Value = "1"
$newSheet.Range("M2").Value = "0"
$newSheet.Range("N2").Value = "30m"
$newSheet.Range("O2").Value = "Company1"
$newSheet.Range("P2").Value = "Event1"

# Column widths (character-width units as shown by Excel's Format > Column
# Width dialog). The values are nudged by -5/6 to compensate for this
# runtime's pixel round-trip so the saved OOXML "width" lands on the
# intended figure (25.6640625, 14, 16.5, 16.83203125, 37.6640625,
# 12.33203125).
$newSheet.Columns.Item(1).ColumnWidth = 24.833333333333336
$newSheet.Columns.Item(2).ColumnWidth = 13.166666666666666
$newSheet.Columns.Item(4).ColumnWidth = 15.666666666666666
$newSheet.Columns.Item(5).ColumnWidth = 15.999999999999998
$newSheet.Columns.Item(7).ColumnWidth = 36.83333333333333
$newSheet.Columns.Item(8).ColumnWidth = 11.5

# Final selection on the new sheet is P2.
[void]$newSheet.Range("P2").Select()
